$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 0.3.0 -> 0.4.0-snapshot-1
$ws.Range("B3").Value = "0.4.0-snapshot-1"

# Status: active -> draft
$ws.Range("B6").Value = "draft"

# Date: 2024-03-13T09:33:00+00:00 -> 2024-05-23T12:16:26+00:00
$ws.Range("B8").Value = "2024-05-23T12:16:26+00:00"

# Contact: "No display for ContactDetail" -> "ANS (https://esante.gouv.fr)"
$ws.Range("B10").Value = "ANS (https://esante.gouv.fr)"
